$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (outside the used range) used to stage text values so that
# numeric/date-looking strings are pasted back as TEXT (matching the source
# data, which stores these as text) instead of being auto-converted by Excels
# smart input parsing. PasteSpecial(values-only) avoids carrying the scratch
# cells "Text" number format onto the destination cells.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$scratch.Value = "0.9992"
$scratch.Copy()
$ws.Range("H2").PasteSpecial(-4163)
$scratch.Value = "0.5883"
$scratch.Copy()
$ws.Range("I2").PasteSpecial(-4163)
$ws.Range("N2").Value = "R$ -710.32"
$ws.Range("O2").Value = "R$ -17.93"
$scratch.Value = "2025-04-04"
$scratch.Copy()
$ws.Range("P2").PasteSpecial(-4163)
$scratch.Value = "0.9976"
$scratch.Copy()
$ws.Range("H3").PasteSpecial(-4163)
$scratch.Value = "0.4645"
$scratch.Copy()
$ws.Range("I3").PasteSpecial(-4163)
$ws.Range("M3").Value = 12
$ws.Range("N3").Value = "R$ -935.91"
$ws.Range("O3").Value = "R$ -99.83"
$scratch.Value = "2025-04-02"
$scratch.Copy()
$ws.Range("P3").PasteSpecial(-4163)
$scratch.Value = "0.9952"
$scratch.Copy()
$ws.Range("H4").PasteSpecial(-4163)
$scratch.Value = "0.3847"
$scratch.Copy()
$ws.Range("I4").PasteSpecial(-4163)
$ws.Range("N4").Value = "R$ -960.01"
$ws.Range("O4").Value = "R$ -64.09"
$scratch.Value = "0.9911"
$scratch.Copy()
$ws.Range("H5").PasteSpecial(-4163)
$scratch.Value = "0.3107"
$scratch.Copy()
$ws.Range("I5").PasteSpecial(-4163)
$ws.Range("J5").Value = "ITM → ATM"
$ws.Range("N5").Value = "R$ -966.68"
$ws.Range("O5").Value = "R$ -69.38"
$scratch.Value = "2025-04-02"
$scratch.Copy()
$ws.Range("P5").PasteSpecial(-4163)
$scratch.Value = "0.9843"
$scratch.Copy()
$ws.Range("H6").PasteSpecial(-4163)
$scratch.Value = "0.2447"
$scratch.Copy()
$ws.Range("I6").PasteSpecial(-4163)
$ws.Range("N6").Value = "R$ -964.62"
$ws.Range("O6").Value = "R$ -108.90"
$scratch.Value = "0.6628"
$scratch.Copy()
$ws.Range("H7").PasteSpecial(-4163)
$ws.Range("L7").Value = 252
$ws.Range("M7").Value = 2
$ws.Range("N7").Value = "R$ 1051.01"
$ws.Range("O7").Value = "R$ 1051.01"
$scratch.Value = "2025-06-18"
$scratch.Copy()
$ws.Range("P7").PasteSpecial(-4163)
$scratch.Value = "0.4271"
$scratch.Copy()
$ws.Range("H8").PasteSpecial(-4163)
$ws.Range("N8").Value = "R$ 437.61"
$ws.Range("O8").Value = "R$ 960.34"
$scratch.Value = "2025-06-10"
$scratch.Copy()
$ws.Range("P8").PasteSpecial(-4163)
$scratch.Value = "0.3690"
$scratch.Copy()
$ws.Range("H9").PasteSpecial(-4163)
$scratch.Value = "0.9982"
$scratch.Copy()
$ws.Range("I9").PasteSpecial(-4163)
$ws.Range("N9").Value = "R$ 278.37"
$ws.Range("O9").Value = "R$ 875.59"
$scratch.Value = "2025-06-10"
$scratch.Copy()
$ws.Range("P9").PasteSpecial(-4163)
$scratch.Value = "0.5519"
$scratch.Copy()
$ws.Range("H10").PasteSpecial(-4163)
$ws.Range("L10").Value = 120
$ws.Range("N10").Value = "R$ 2020.20"
$ws.Range("O10").Value = "R$ 2020.20"
$scratch.Value = "0.2569"
$scratch.Copy()
$ws.Range("H11").PasteSpecial(-4163)
$scratch.Value = "0.7320"
$scratch.Copy()
$ws.Range("I11").PasteSpecial(-4163)
$ws.Range("N11").Value = "R$ 590.57"
$ws.Range("O11").Value = "R$ 724.28"
$scratch.Value = "0.2016"
$scratch.Copy()
$ws.Range("H12").PasteSpecial(-4163)
$scratch.Value = "0.1700"
$scratch.Copy()
$ws.Range("I12").PasteSpecial(-4163)
$ws.Range("K12").Value = 0.15
$ws.Range("M12").Value = 3
$ws.Range("N12").Value = "R$ 578.31"
$ws.Range("O12").Value = "R$ 578.31"
$scratch.Value = "2025-06-18"
$scratch.Copy()
$ws.Range("P12").PasteSpecial(-4163)

# Clean up the scratch cell so it leaves no trace in the saved sheet.
$scratch.Clear()
